$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated PSSM values for range B2:K21 (rows follow residues F,W,Y,P,M,I,L,V,A,G,C,S,T,N,Q,D,E,H,K,R
# in column A; columns B:K are sequence positions 0-9). Values refreshed per supplemental figures.
$newValues = @(
    @(-17.08291826251373, 1.796297015609628, -17.08291826251373, -17.08291826251373, -17.08291826251373, -17.08291826251373, -17.08291826251373, -17.08291826251373, -17.08291826251373, -17.08291826251373),
    @(-17.08291826251373, -17.08291826251373, -17.08291826251373, -17.08291826251373, -17.08291826251373, -17.08291826251373, -17.08291826251373, 2.249234139597022, -17.08291826251373, -17.08291826251373),
    @(-17.08291826251373, 2.013387679208233, 2.206688663983783, -17.08291826251373, 3.560480848803758, -17.08291826251373, 4.321918222349233, -17.08291826251373, 2.097090417437743, -17.08291826251373),
    @(-17.08291826251373, 1.929228841449173, -17.08291826251373, -17.08291826251373, -17.08291826251373, 2.705247242947324, -17.08291826251373, -17.08291826251373, -17.08291826251373, -17.08291826251373),
    @(-17.08291826251373, -17.08291826251373, -17.08291826251373, -17.08291826251373, -17.08291826251373, -17.08291826251373, -17.08291826251373, -17.08291826251373, -17.08291826251373, -17.08291826251373),
    @(2.543968923154591, -17.08291826251373, -17.08291826251373, -17.08291826251373, -17.08291826251373, -17.08291826251373, -17.08291826251373, -17.08291826251373, -17.08291826251373, -17.08291826251373),
    @(-17.08291826251373, -17.08291826251373, -17.08291826251373, 1.865463228938493, -17.08291826251373, -17.08291826251373, -17.08291826251373, -17.08291826251373, -17.08291826251373, -17.08291826251373),
    @(3.824559981303727, -17.08291826251373, -17.08291826251373, -17.08291826251373, -17.08291826251373, -17.08291826251373, -17.08291826251373, -17.08291826251373, -17.08291826251373, -17.08291826251373),
    @(-17.08291826251373, -17.08291826251373, -17.08291826251373, -17.08291826251373, -17.08291826251373, -17.08291826251373, -17.08291826251373, 1.122612512532007, -17.08291826251373, 1.574264090510578),
    @(-17.08291826251373, -17.08291826251373, -17.08291826251373, 2.922965289377769, -17.08291826251373, 2.973800146863287, -17.08291826251373, -17.08291826251373, -17.08291826251373, 1.519929184323978),
    @(-17.08291826251373, -17.08291826251373, -17.08291826251373, -17.08291826251373, -17.08291826251373, -17.08291826251373, -17.08291826251373, -17.08291826251373, -17.08291826251373, -17.08291826251373),
    @(-17.08291826251373, -17.08291826251373, -17.08291826251373, 2.352509874217248, -17.08291826251373, -17.08291826251373, -17.08291826251373, -17.08291826251373, 2.048026780498276, 1.891521535693599),
    @(-17.08291826251373, -17.08291826251373, 1.096509784188927, -17.08291826251373, -17.08291826251373, -17.08291826251373, -17.08291826251373, -17.08291826251373, -17.08291826251373, 2.451584110661615),
    @(-17.08291826251373, -17.08291826251373, 0.9943199371982013, -17.08291826251373, -17.08291826251373, -17.08291826251373, -17.08291826251373, -17.08291826251373, -17.08291826251373, -17.08291826251373),
    @(-17.08291826251373, -17.08291826251373, -17.08291826251373, -17.08291826251373, -17.08291826251373, -17.08291826251373, -17.08291826251373, -17.08291826251373, 1.93386010652161, -17.08291826251373),
    @(-17.08291826251373, 2.219954898151224, 2.564492731976138, -17.08291826251373, -17.08291826251373, -17.08291826251373, -17.08291826251373, 2.195433636980447, 2.420629941752761, -17.08291826251373),
    @(-17.08291826251373, -17.08291826251373, -17.08291826251373, -17.08291826251373, -17.08291826251373, -17.08291826251373, -17.08291826251373, 1.449228360616537, 1.269767479033787, -17.08291826251373),
    @(-17.08291826251373, -17.08291826251373, 1.464898832928643, -17.08291826251373, -17.08291826251373, -17.08291826251373, -17.08291826251373, 1.902681893435453, -17.08291826251373, -17.08291826251373),
    @(-17.08291826251373, 0.9238734869235141, 1.365879631434274, -17.08291826251373, 3.035937920045323, -17.08291826251373, -17.08291826251373, 1.013081337088041, -17.08291826251373, 2.314365715753257),
    @(-17.08291826251373, 1.087344244063502, -17.08291826251373, 1.873767007314233, -17.08291826251373, 2.491122692549226, -17.08291826251373, -17.08291826251373, -17.08291826251373, -17.08291826251373)
)

$startRow = 2
$startCol = 2  # column B

for ($i = 0; $i -lt $newValues.Count; $i++) {
    $rowValues = $newValues[$i]
    $row = $startRow + $i
    for ($j = 0; $j -lt $rowValues.Count; $j++) {
        $col = $startCol + $j
        $ws.Cells.Item($row, $col).Value = $rowValues[$j]
    }
}
